# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (holding-detail rows, same layout as the
# "2021-Q2" sheet) between the existing "2021-Q2" and "总计" sheets, and
# updates the "总计" (totals) sheet with a new leading row summarising the
# 2022-Q1 quarter (the 2021-Q2 totals row stays, now second).

$wb = $excel.ActiveWorkbook

# --- Step 1: all sheet-structure changes first -----------------------
# Duplicate the "总计" sheet — this keeps the bold/bordered header style
# (style index 2) the new sheet needs, since it is the "总计" sheet's own
# header style, not "2021-Q2"'s. The duplicate (fresh copy) lands right
# after "总计"; move the original "总计" sheet up to sit right after
# "2021-Q2" instead, then rename in place, so that the ORIGINAL "总计"
# worksheet (lower internal sheetId) becomes "2022-Q1" and the NEW copy
# (higher sheetId) becomes "总计" — matching how sheetIds land when a
# sheet is inserted ahead of an existing one in Excel.
# (Worksheet handles obtained before a Copy/Add/Delete/Move resolve by
# position, so re-fetch everything by name only AFTER the sheet
# collection is final.)
$wb.Worksheets.Item("总计").Copy($null, $wb.Worksheets.Item("2021-Q2"))
$wb.Worksheets.Item(3).Move($wb.Worksheets.Item(2))
$wb.Worksheets.Item(2).Name = "2022-Q1"
$wb.Worksheets.Item(3).Name = "总计"

# --- Step 2: fetch stable handles now that sheet order is final -------
$sheetNew   = $wb.Worksheets.Item("2022-Q1")
$sheetTotal = $wb.Worksheets.Item("总计")

# --- helper: write a value as literal text (not auto-parsed as a
# number), without leaving a residual NumberFormat style on the cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- "2022-Q1" sheet: header row ----
# Extend the existing header style (currently only on B1:D1, copied from
# "总计") across E1:H1 before filling in the remaining header labels.
$sheetNew.Range("B1").Copy()
$sheetNew.Range("E1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$sheetNew.Range("B1").Value = "基金代码"
$sheetNew.Range("C1").Value = "基金名称"
$sheetNew.Range("D1").Value = "基金规模"
$sheetNew.Range("E1").Value = "股票总仓位"
$sheetNew.Range("F1").Value = "仓位占比"
$sheetNew.Range("G1").Value = "持有市值(亿元)"
$sheetNew.Range("H1").Value = "仓位排名"

# ---- "2022-Q1" sheet: data rows ----
# Row 3 is brand new (only row 2 existed on the template "总计" sheet), so
# give A3 the same style as A2 before filling in values.
$sheetNew.Range("A2").Copy()
$sheetNew.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$sheetNew.Range("A2").Value = 0
Set-TextValue $sheetNew.Range("B2") "003720"
Set-TextValue $sheetNew.Range("C2") "易方达标普生物科技指数（QDII-LOF）美元"
Set-TextValue $sheetNew.Range("D2") "2.11"
Set-TextValue $sheetNew.Range("E2") "94.00"
Set-TextValue $sheetNew.Range("F2") "1.09"
Set-TextValue $sheetNew.Range("G2") "0.0230"
$sheetNew.Range("H2").Value = 1

$sheetNew.Range("A3").Value = 1
Set-TextValue $sheetNew.Range("B3") "161127"
Set-TextValue $sheetNew.Range("C3") "易方达标普生物科技指数（QDII-LOF）人民币"
Set-TextValue $sheetNew.Range("D3") "2.11"
Set-TextValue $sheetNew.Range("E3") "94.00"
Set-TextValue $sheetNew.Range("F3") "1.09"
Set-TextValue $sheetNew.Range("G3") "0.0230"
$sheetNew.Range("H3").Value = 1

# ---- "总计" sheet: insert the new 2022-Q1 totals row above 2021-Q2 ----
# Push the existing 2021-Q2 summary row down to row 3 first, then write the
# new 2022-Q1 summary row into row 2. Row 3 is brand new, so copy A2's
# style onto A3 before filling in values.
$sheetTotal.Range("A2").Copy()
$sheetTotal.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("B3").Value = "2021-Q2"
$sheetTotal.Range("C3").Value = 2
$sheetTotal.Range("D3").Value = 0.03

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.05
